$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J holds "municipio-nombre" metadata. The new curated dimensions
# reclassify it as a proper SDMX dimension (refArea) rather than a measure,
# matching the treatment already given to "provincia-nombre" (column L)
# and "comarca-nombre" (column N).
$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "dim"
$ws.Range("J4").Value = "URI-Municipio"
